$wb = $excel.ActiveWorkbook

# --- 1. Append the new log entry as row 28 on the "Logs" sheet ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A28").Value = "Opvolging klacht"
$ws.Range("B28").Value = "mailmind.test@zohomail.eu"
$ws.Range("D28").Value = "Klacht / Probleem"
$ws.Range("F28").Value = "2025-08-28 21:20:37"
$ws.Range("G28").Value = "Nee"
$ws.Range("H28").Value = "Ja"
$ws.Range("I28").Value = "Nee"
$ws.Range("J28").Value = "Nee"

# --- 2. Extend the conditional formatting ranges to include the new row ---
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $ws.Range($col + "2:" + $col + "27")
    $newRange = $ws.Range($col + "2:" + $col + "28")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- 3. Update the "Dashboard" summary table: the new complaint pushes
#        "Klacht / Probleem" to 2 occurrences, ahead of "Overig" (1) ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Klacht / Probleem"
$dash.Range("B6").Value = 2
$dash.Range("A7").Value = "Overig"
$dash.Range("B7").Value = 1
